$wb = $excel.ActiveWorkbook

# --- Sheet "Sheet1": data changes ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A3").Value = "video_file"
$ws1.Range("B6").Value = "Cam_Basic.sdf"
$ws1.Range("C6").Value = "1,1,1,0,0,0"
$ws1.Range("A7").Value = "markers"
$ws1.Range("C7").Value = "1,1,0.5,0.72,0,0"
$ws1.Range("C8").Value = "-1,0,0.1,0,0.52,0"

# --- Selection changes (no data change) ---
$wsMain = $wb.Worksheets.Item("Main")
$wsMain.Activate()
$wsMain.Range("A3").Select()

$ws1.Activate()
$ws1.Range("C6").Select()
